$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "def_value"
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = "def_pixel_size"
$ws.Range("B6").Value = 30.716988260000001

$ws.Range("B7").Value = "uint8"
$ws.Range("A7").Value = "def_dtype"

$null = $ws.Range("B6").Select()

$null = $ws.Hyperlinks.Add($ws.Range("B2"), "http://bart.ideam.gov.co/cneideam/Capasgeo/")
